# Append two new match rows (162, 163) to the Romania Liga 1 2023-2024 sheet,
# mirroring the formatting of the last existing data row (161).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use the last existing row as a formatting template for the new rows so the
# cell styles (bold/bordered index column, datetime-formatted match-date
# column) carry over exactly as they appear elsewhere in the sheet.
$template = $ws.Range("A161:V161")
$template.Copy($ws.Range("A162:V162"))
$template.Copy($ws.Range("A163:V163"))

# --- Row 162: U Craiova 1948 2-0 Petrolul ---
$ws.Range("A162").Value = 161
$ws.Range("B162").Value = "romania"
$ws.Range("C162").Value = "liga-1"
$ws.Range("D162").Value = "2023-2024"
$ws.Range("E162").Value = 45279.6875
$ws.Range("F162").Value = "U Craiova 1948"
$ws.Range("G162").Value = 2
$ws.Range("H162").Value = "Petrolul"
$ws.Range("I162").Value = 0
$ws.Range("J162").Value = 2.17
$ws.Range("K162").Value = "15/12/2023 19:12"
$ws.Range("L162").Value = 1.94
$ws.Range("M162").Value = "19/12/2023 16:28"
$ws.Range("N162").Value = 3.26
$ws.Range("O162").Value = "15/12/2023 19:12"
$ws.Range("P162").Value = 3.39
$ws.Range("Q162").Value = "19/12/2023 16:28"
$ws.Range("R162").Value = 3.25
$ws.Range("S162").Value = "15/12/2023 19:12"
$ws.Range("T162").Value = 4.23
$ws.Range("U162").Value = "19/12/2023 16:28"
$ws.Range("V162").Value = "https://www.betexplorer.com/football/romania/liga-1/fc-u-craiova-petrolul/4Ex4Fjfg/"

# --- Row 163: Poli Iasi 1-3 FCSB ---
$ws.Range("A163").Value = 162
$ws.Range("B163").Value = "romania"
$ws.Range("C163").Value = "liga-1"
$ws.Range("D163").Value = "2023-2024"
$ws.Range("E163").Value = 45279.8125
$ws.Range("F163").Value = "Poli Iasi"
$ws.Range("G163").Value = 1
$ws.Range("H163").Value = "FCSB"
$ws.Range("I163").Value = 3
$ws.Range("J163").Value = 4.5
$ws.Range("K163").Value = "16/12/2023 19:13"
$ws.Range("L163").Value = 4.68
$ws.Range("M163").Value = "19/12/2023 19:28"
$ws.Range("N163").Value = 3.67
$ws.Range("O163").Value = "16/12/2023 19:13"
$ws.Range("P163").Value = 3.68
$ws.Range("Q163").Value = "19/12/2023 19:25"
$ws.Range("R163").Value = 1.71
$ws.Range("S163").Value = "16/12/2023 19:13"
$ws.Range("T163").Value = 1.78
$ws.Range("U163").Value = "19/12/2023 19:25"
$ws.Range("V163").Value = "https://www.betexplorer.com/football/romania/liga-1/poli-iasi-fcsb/MkmaGWvm/"
